$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.415.10"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.850.12"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.09%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "233.09"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("E6").Value = "  +0.10%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4756"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +2.75%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2742"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +1.94%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06320"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +1.62%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "17.56"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +9.46%  "
$ws.Range("D11").Value = "1.865.96"
$ws.Range("E11").Value = "  +1.90%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07454"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.44%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.946"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.11%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "84.49"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +1.92%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.6225"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "30.392.63"
$ws.Range("E16").Value = "  +1.21%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "244.19"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +7.81%  "
$ws.Range("E18").Value = "  +0.13%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.64"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +2.97%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "0.000007316"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.75%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.900"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +1.73%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.894"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.41%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "164.86"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "9.068"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.50%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "17.96"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("E27").Value = "  +1.32%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.1028"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.03%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.344"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.64%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.030"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -0.34%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "3.806"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +1.48%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.04825"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  +0.17%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.6952"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  +0.66%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.01894"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +5.18%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.682"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.82%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "1.996"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +4.52%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.8734"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.98%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "106.55"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +3.68%  "
$ws.Range("E41").Value = "  +0.11%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.508"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.72%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.4049"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "7.135"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +3.61%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "63.06"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +6.82%  "
$ws.Range("E46").Value = "  +0.69%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "33.69"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +3.90%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "8.522"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.05510"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -0.11%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.3672"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +1.42%  "
